$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 223; this shifts existing rows 223:242 down to 224:243
$ws.Rows(223).Insert()

# Populate the newly inserted row 223 with the new record
$ws.Range("A223").Value = 7
$ws.Range("B223").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C223").Value = "Ñuble"
$ws.Range("D223").Value = 44769
$ws.Range("E223").Value = 16
$ws.Range("F223").Value = 100112043
$ws.Range("G223").Value = "Pepino ensalada"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 80
$ws.Range("K223").Value = 19000
$ws.Range("L223").Value = 20000
$ws.Range("M223").Value = 19500
$ws.Range("N223").Value = "`$/caja 60 unidades"
$ws.Range("O223").Value = "Región de Arica y Parinacota"
$ws.Range("P223").Value = 325
$ws.Range("Q223").Value = 60
$ws.Range("R223").Value = "Hortaliza"

# Make sure the D223 style matches date-style cells (style index "2" with YYYY-MM-DD HH:MM:SS numFmt)
$ws.Range("D223").NumberFormat = $ws.Range("D224").NumberFormat
